# Daily attendance processing - 2025-10-18 11:39:58
# Reorders the "Recorded By" (column G) values in specific rows of the
# "Session Analysis Results" sheet so that the names appear in the
# updated order captured in the source export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    "G2" = 'System, backup@backdoor.com, system'
    "G3" = 'dnasr281@gmail.com, System'
    "G4" = 'System, backup@backdoor.com'
    "G6" = 'dnasr281@gmail.com, System'
    "G10" = 'dnasr281@gmail.com, System'
    "G11" = 'dnasr281@gmail.com, System'
    "G12" = 'dnasr281@gmail.com, System'
    "G13" = 'dnasr281@gmail.com, System'
    "G14" = 'dnasr281@gmail.com, System'
    "G15" = 'dnasr281@gmail.com, System'
    "G17" = 'dnasr281@gmail.com, System'
    "G18" = 'dnasr281@gmail.com, System'
    "G19" = 'dnasr281@gmail.com, System'
    "G20" = 'dnasr281@gmail.com, System'
    "G29" = 'System, backup@backdoor.com, system'
    "G30" = 'dnasr281@gmail.com, System'
    "G31" = 'System, backup@backdoor.com'
    "G33" = 'dnasr281@gmail.com, System'
    "G37" = 'dnasr281@gmail.com, System'
    "G38" = 'dnasr281@gmail.com, System'
    "G39" = 'dnasr281@gmail.com, System'
    "G40" = 'dnasr281@gmail.com, System'
    "G41" = 'dnasr281@gmail.com, System'
    "G42" = 'dnasr281@gmail.com, System'
    "G44" = 'dnasr281@gmail.com, System'
    "G45" = 'dnasr281@gmail.com, System'
    "G46" = 'dnasr281@gmail.com, System'
    "G47" = 'dnasr281@gmail.com, System'
    "G56" = 'System, backup@backdoor.com, system'
    "G57" = 'dnasr281@gmail.com, System'
    "G58" = 'System, backup@backdoor.com'
    "G60" = 'dnasr281@gmail.com, System'
    "G64" = 'dnasr281@gmail.com, System'
    "G65" = 'dnasr281@gmail.com, System'
    "G66" = 'dnasr281@gmail.com, System'
    "G67" = 'dnasr281@gmail.com, System'
    "G68" = 'dnasr281@gmail.com, System'
    "G69" = 'dnasr281@gmail.com, System'
    "G71" = 'dnasr281@gmail.com, System'
    "G72" = 'dnasr281@gmail.com, System'
    "G73" = 'dnasr281@gmail.com, System'
    "G74" = 'dnasr281@gmail.com, System'
    "G83" = 'System, backup@backdoor.com'
    "G84" = 'System, backup@backdoor.com'
    "G86" = 'dnasr281@gmail.com, System'
    "G87" = 'dnasr281@gmail.com, System'
    "G88" = 'dnasr281@gmail.com, System'
    "G89" = 'dnasr281@gmail.com, System'
    "G93" = 'dnasr281@gmail.com, System'
    "G95" = 'dnasr281@gmail.com, System'
    "G96" = 'dnasr281@gmail.com, System'
    "G97" = 'dnasr281@gmail.com, System'
    "G99" = 'dnasr281@gmail.com, System'
    "G109" = 'System, backup@backdoor.com'
    "G110" = 'System, backup@backdoor.com'
    "G112" = 'dnasr281@gmail.com, System'
    "G113" = 'dnasr281@gmail.com, System'
    "G114" = 'dnasr281@gmail.com, System'
    "G115" = 'dnasr281@gmail.com, System'
    "G119" = 'dnasr281@gmail.com, System'
    "G121" = 'dnasr281@gmail.com, System'
    "G122" = 'dnasr281@gmail.com, System'
    "G123" = 'dnasr281@gmail.com, System'
    "G125" = 'dnasr281@gmail.com, System'
    "G135" = 'System, backup@backdoor.com'
    "G136" = 'System, backup@backdoor.com'
    "G138" = 'dnasr281@gmail.com, System'
    "G139" = 'dnasr281@gmail.com, System'
    "G140" = 'dnasr281@gmail.com, System'
    "G141" = 'dnasr281@gmail.com, System'
    "G145" = 'dnasr281@gmail.com, System'
    "G147" = 'dnasr281@gmail.com, System'
    "G148" = 'dnasr281@gmail.com, System'
    "G149" = 'dnasr281@gmail.com, System'
    "G151" = 'dnasr281@gmail.com, System'
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
